$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("POP")
$ws.Range("B2").Value = 1.995
$ws.Range("C2").Value = 0.3286
$ws.Range("B3").Value = 0.5973
$ws.Range("C3").Value = 0.1661
$ws.Range("D3").Value = 0.0003
$ws.Range("B4").Value = -0.1298
$ws.Range("C4").Value = 0.233
$ws.Range("D4").Value = 0.5774
$ws.Range("B5").Value = 0.3552
$ws.Range("C5").Value = 0.1763
$ws.Range("D5").Value = 0.0439
$ws.Range("B6").Value = -0.2309
$ws.Range("C6").Value = 0.3212
$ws.Range("D6").Value = 0.4723
$ws.Range("B7").Value = 0.5263
$ws.Range("C7").Value = 0.1924
$ws.Range("D7").Value = 0.0062
$ws.Range("B8").Value = 0.1948
$ws.Range("C8").Value = 0.1771
$ws.Range("D8").Value = 0.2714
$ws.Range("B9").Value = -0.1604
$ws.Range("C9").Value = 0.2829
$ws.Range("D9").Value = 0.5707
$ws.Range("B10").Value = 0.4122
$ws.Range("C10").Value = 0.3443
$ws.Range("D10").Value = 0.2311
$ws.Range("B11").Value = -0.0234
$ws.Range("C11").Value = 0.2476
$ws.Range("D11").Value = 0.9248
$ws.Range("B12").Value = -0.0139
$ws.Range("C12").Value = 0.3245
$ws.Range("D12").Value = 0.9659
$ws.Range("B13").Value = 0.1324
$ws.Range("C13").Value = 0.315
$ws.Range("D13").Value = 0.6743
$ws.Range("B14").Value = 0.2194
$ws.Range("C14").Value = 0.2933
$ws.Range("D14").Value = 0.4545
$ws.Range("B15").Value = 0.0665
$ws.Range("C15").Value = 0.2919
$ws.Range("D15").Value = 0.8197
$ws.Range("B16").Value = -0.0084
$ws.Range("C16").Value = 0.308
$ws.Range("D16").Value = 0.9784
$ws.Range("B17").Value = -0.0668
$ws.Range("C17").Value = 0.3466
$ws.Range("D17").Value = 0.8473
$ws.Range("B18").Value = -0.1587
$ws.Range("C18").Value = 0.3993
$ws.Range("D18").Value = 0.691
$ws.Range("B19").Value = -0.1932
$ws.Range("C19").Value = 0.4323
$ws.Range("D19").Value = 0.6549
$ws.Range("B20").Value = -0.1503
$ws.Range("C20").Value = 0.4522
$ws.Range("D20").Value = 0.7395
$ws.Range("B21").Value = -0.1439
$ws.Range("C21").Value = 0.5112
$ws.Range("D21").Value = 0.7783
$ws.Range("B22").Value = -0.1448
$ws.Range("C22").Value = 0.6737
$ws.Range("D22").Value = 0.8298
$ws.Range("B23").Value = -0.886
$ws.Range("C23").Value = 0.0543
$ws = $wb.Worksheets.Item("Pesticide")
$ws.Range("B2").Value = -1.4781
$ws.Range("C2").Value = 6.7583
$ws.Range("D2").Value = 0.8269
$ws.Range("B3").Value = -1.1122
$ws.Range("C3").Value = 25.8366
$ws.Range("D3").Value = 0.9657
$ws.Range("B4").Value = 0.6574
$ws.Range("C4").Value = 0.3803
$ws.Range("D4").Value = 0.0839
$ws.Range("B5").Value = -1.1195
$ws.Range("C5").Value = 30.2141
$ws.Range("D5").Value = 0.9704
$ws.Range("B6").Value = 0.3754
$ws.Range("C6").Value = 0.3589
$ws.Range("D6").Value = 0.2956
$ws.Range("B7").Value = 0.1746
$ws.Range("C7").Value = 44.7064
$ws.Range("D7").Value = 0.9969
$ws.Range("B8").Value = -1.5031
$ws.Range("C8").Value = 26.1074
$ws.Range("D8").Value = 0.9541
$ws.Range("B9").Value = -1.3778
$ws.Range("C9").Value = 46.1469
$ws.Range("D9").Value = 0.9762
$ws.Range("B10").Value = -0.0669
$ws.Range("C10").Value = 1.2079
$ws.Range("D10").Value = 0.9558
$ws.Range("B11").Value = 0.5245
$ws.Range("C11").Value = 3.1689
$ws.Range("D11").Value = 0.8685
$ws.Range("B12").Value = 1.0491
$ws.Range("C12").Value = 5.4201
$ws.Range("D12").Value = 0.8465
$ws.Range("B13").Value = 1.5708
$ws.Range("C13").Value = 6.6226
$ws.Range("D13").Value = 0.8125
$ws.Range("B14").Value = 2.0869
$ws.Range("C14").Value = 6.8226
$ws.Range("D14").Value = 0.7597
$ws.Range("B15").Value = 2.3804
$ws.Range("C15").Value = 6.7545
$ws.Range("D15").Value = 0.7245
$ws.Range("B16").Value = 2.6176
$ws.Range("C16").Value = 6.7715
$ws.Range("D16").Value = 0.6991
$ws.Range("B17").Value = 2.8467
$ws.Range("C17").Value = 6.808
$ws.Range("D17").Value = 0.6758
$ws.Range("B18").Value = 3.1535
$ws.Range("C18").Value = 6.8716
$ws.Range("D18").Value = 0.6463
$ws.Range("B19").Value = 2.4883
$ws.Range("C19").Value = 6.9017
$ws.Range("D19").Value = 0.7184
$ws.Range("B20").Value = 1.3112
$ws.Range("C20").Value = 7.0847
$ws.Range("D20").Value = 0.8532
$ws.Range("B21").Value = 0.1292
$ws.Range("C21").Value = 7.8209
$ws.Range("D21").Value = 0.9868
$ws.Range("B22").Value = -1.0528
$ws.Range("C22").Value = 9.3683
$ws.Range("D22").Value = 0.9105
$ws.Range("B23").Value = -0.0885
$ws.Range("C23").Value = 0.1748
$ws.Range("D23").Value = 0.6126
